$wb = $excel.ActiveWorkbook

# Row 91 data for each of the 4 worksheets (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2)
$sheetsData = @(
    @{
        Sheet = "FE_LFT_#1"
        A = 45877.4919212963
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x0C"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 268
        I = 15
    },
    @{
        Sheet = "FE_LFT_#2"
        A = 45877.4919212963
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x1C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 284
        I = 14
    },
    @{
        Sheet = "FE_PLT_#1"
        A = 45877.4919212963
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x5D"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 93
        I = 3
    },
    @{
        Sheet = "FE_PLT_#2"
        A = 45877.4919212963
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5C"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 92
        I = 3
    }
)

foreach ($row in $sheetsData) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    $ws.Cells.Item(91, 1).Value = $row.A
    $ws.Cells.Item(91, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item(91, 2).Value = $row.B
    $ws.Cells.Item(91, 3).Value = $row.C
    $ws.Cells.Item(91, 4).Value = $row.D
    $ws.Cells.Item(91, 5).Value = $row.E
    $ws.Cells.Item(91, 6).Value = $row.F
    $ws.Cells.Item(91, 7).Value = $row.G
    $ws.Cells.Item(91, 8).Value = $row.H
    $ws.Cells.Item(91, 9).Value = $row.I
}
